# Update summary.xlsx with latest data revisions
#
# - B1 label tweak ("Min-Max Allowance" -> "Min - Max Allowance")
# - New columns H:Q (header row + per-row calibration data out to row 4)
# - Revised Test 1-5 readings for rows 2-4
#
# NOTE: every value here is logically *text* (numbers-as-strings, ISO dates
# as strings, etc.) mirroring the original inline-string sheet, so each
# write forces NumberFormat "@" first to stop Excel's literal-to-number/
# date autoconversion; the style is restored immediately after so plain
# data cells keep their original (default) styling instead of inheriting
# the "@" text format as a new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

function Set-HeaderTextValue {
    # Like Set-TextValue, but (re)applies the bold/border/center header
    # styling (matching A1) instead of trying to preserve/restore whatever
    # style the cell already had - needed both for the brand-new H1:Q1
    # cells and for B1, whose $Range.Style round-trip isn't reliable here.
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $ws.Range("A1").Copy()
    $Range.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# --- B1 header label edit ---
Set-HeaderTextValue $ws.Range("B1") "Min - Max Allowance"

# --- New header cells H1:Q1, styled like the existing header row ---
$newHeaders = @("Manufacturer", "Serial Number", "Model", "Calibration Date", "Calibration Due", "Unit Number", "Customer/Company", "Phone Number", "Address", "OCR Text")
$headerCols = @("H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $cell = $ws.Range($headerCols[$i] + "1")
    Set-HeaderTextValue $cell $newHeaders[$i]
}

# --- Revised Test 1-5 readings (rows 2-4, columns C-G) ---
Set-TextValue $ws.Range("C2") "567.8"
Set-TextValue $ws.Range("D2") "553.3"
Set-TextValue $ws.Range("E2") "550.2"
Set-TextValue $ws.Range("F2") "561.0"
Set-TextValue $ws.Range("G2") "559.0"

Set-TextValue $ws.Range("C3") "362.5"
Set-TextValue $ws.Range("D3") "353.5"
Set-TextValue $ws.Range("E3") "351.1"
Set-TextValue $ws.Range("F3") "360.8"
Set-TextValue $ws.Range("G3") "353.5"

Set-TextValue $ws.Range("C4") "203.3"
Set-TextValue $ws.Range("D4") "200.2"
Set-TextValue $ws.Range("E4") "198.3"
Set-TextValue $ws.Range("F4") "204.6"
Set-TextValue $ws.Range("G4") "197.7"

# --- New columns H:Q for rows 2-4 ---
# H/I/J (Manufacturer/Serial Number/Model) and M:Q stay blank; only
# Calibration Date (K) / Calibration Due (L) are populated, same value
# repeated for every test row.
for ($row = 2; $row -le 4; $row++) {
    Set-TextValue $ws.Cells.Item($row, 8) ""    # H - Manufacturer
    Set-TextValue $ws.Cells.Item($row, 9) ""    # I - Serial Number
    Set-TextValue $ws.Cells.Item($row, 10) ""   # J - Model
    Set-TextValue $ws.Cells.Item($row, 11) "2025-03-03"  # K - Calibration Date
    Set-TextValue $ws.Cells.Item($row, 12) "2026-03-03"  # L - Calibration Due
    Set-TextValue $ws.Cells.Item($row, 13) ""   # M - Unit Number
    Set-TextValue $ws.Cells.Item($row, 14) ""   # N - Customer/Company
    Set-TextValue $ws.Cells.Item($row, 15) ""   # O - Phone Number
    Set-TextValue $ws.Cells.Item($row, 16) ""   # P - Address
    Set-TextValue $ws.Cells.Item($row, 17) ""   # Q - OCR Text
}
